$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.3880251944065094
$ws.Range("B1").Value = 0.7843900918960571
$ws.Range("C1").Value = 3.191421747207642
$ws.Range("D1").Value = 2.629640579223633
$ws.Range("E1").Value = 0.9633427858352661
